# Applies the "subject 27 - 29" edit: fills in previously-blank survey
# answers for respondents (rows) 27-29 (worksheet rows 29-31), flips the
# Condition (G/F) for rows 28-29 and 30 (worksheet rows 30, 31, 32), and
# appends a new respondent (subject 31) in worksheet row 33.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 29 (subject 27): fill in previously blank answers ---
$ws.Range("B29").Value = "Female"
$ws.Range("C29").Value = 20
$ws.Range("D29").Value = "Biomedical Eng."
$ws.Range("E29").Value = 1
$ws.Range("H29").Value = 1
$ws.Range("I29").Value = 1
$ws.Range("J29").Value = 2
# K29/L29/A29 already correct ("F" / 1 / 27)

# --- Row 30 (subject 28): fill in previously blank answers + fix Condition ---
$ws.Range("B30").Value = "Male"
$ws.Range("C30").Value = 21
$ws.Range("D30").Value = "Mechanical Eng."
$ws.Range("E30").Value = 4
$ws.Range("F30").Value = "Game Console"
$ws.Range("G30").Value = "Joystick; "
$ws.Range("H30").Value = 4
$ws.Range("I30").Value = 1
$ws.Range("J30").Value = 4
$ws.Range("K30").Value = "F"

# --- Row 31 (subject 29): fill in previously blank answers + fix Condition ---
$ws.Range("B31").Value = "Female"
$ws.Range("C31").Value = 19
$ws.Range("D31").Value = "Mechanical Eng."
$ws.Range("E31").Value = 1
$ws.Range("F31").Value = "Desktop"
$ws.Range("G31").Value = "Keyboard/Mouse"
$ws.Range("H31").Value = 1
$ws.Range("I31").Value = 1
$ws.Range("J31").Value = 2
$ws.Range("K31").Value = "G"

# --- Row 32 (subject 30): fix Condition only ---
$ws.Range("K32").Value = "F"

# --- Row 33 (new subject 31) ---
$ws.Range("A33").Value = 31
$ws.Range("K33").Value = "G"
$ws.Range("L33").Value = 1

# Rows 29-31 now wrap onto two lines in the UI, same as other answered rows.
$ws.Range("A29:L31").EntireRow.RowHeight = 30

# Match the saved view state: scrolled down with J31 selected.
$ws.Range("J31").Select() | Out-Null
